$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# Rows 2-37: only the "Field Name" (B) and "Field Type" (D) columns are
# re-ordered within each class's block of rows; "Class Name" (A) and
# "Field Modifier" (C) stay the same per row.

$fieldNames = @{
    2  = "productAttributeId"
    3  = "name"
    4  = "value"
    5  = "type"
    6  = "serialVersionUID"
    7  = "id"
    8  = "esProductService"
    9  = "pic"
    10 = "id"
    11 = "productCategoryName"
    12 = "productCategoryId"
    13 = "subTitle"
    14 = "price"
    15 = "newStatus"
    16 = "attrValueList"
    17 = "brandId"
    18 = "serialVersionUID"
    19 = "stock"
    20 = "brandName"
    21 = "name"
    22 = "keywords"
    23 = "promotionType"
    24 = "sale"
    25 = "recommandStatus"
    26 = "productSn"
    27 = "sort"
    28 = "elasticsearchRestTemplate"
    29 = "productDao"
    30 = "LOGGER"
    31 = "productRepository"
    32 = "attrName"
    33 = "attrId"
    34 = "attrValues"
    35 = "brandNames"
    36 = "productCategoryNames"
    37 = "productAttrs"
}

$fieldTypes = @{
    2  = "java.lang.Long"
    3  = "java.lang.String"
    4  = "java.lang.String"
    5  = "java.lang.Integer"
    6  = "long"
    7  = "java.lang.Long"
    8  = "com.macro.mall.search.service.EsProductService"
    9  = "java.lang.String"
    10 = "java.lang.Long"
    11 = "java.lang.String"
    12 = "java.lang.Long"
    13 = "java.lang.String"
    14 = "java.math.BigDecimal"
    15 = "java.lang.Integer"
    16 = "java.util.List"
    17 = "java.lang.Long"
    18 = "long"
    19 = "java.lang.Integer"
    20 = "java.lang.String"
    21 = "java.lang.String"
    22 = "java.lang.String"
    23 = "java.lang.Integer"
    24 = "java.lang.Integer"
    25 = "java.lang.Integer"
    26 = "java.lang.String"
    27 = "java.lang.Integer"
    28 = "org.springframework.data.elasticsearch.core.ElasticsearchRestTemplate"
    29 = "com.macro.mall.search.dao.EsProductDao"
    30 = "org.slf4j.Logger"
    31 = "com.macro.mall.search.repository.EsProductRepository"
    32 = "java.lang.String"
    33 = "java.lang.Long"
    34 = "java.util.List"
    35 = "java.util.List"
    36 = "java.util.List"
    37 = "java.util.List"
}

foreach ($r in 2..37) {
    $ws.Cells.Item($r, 2).Value = $fieldNames[$r]
    $ws.Cells.Item($r, 4).Value = $fieldTypes[$r]
}
